{"js": "// Move the \"Langages : ...\" skills line so it becomes the first item of the\n// \"COMPETENCES TECHNIQUES\" list (i.e. insert it right before the\n// \"Visualisation : tableau\" paragraph) instead of sitting between \"MLOps\"\n// and \"Bases de donn\u00e9es\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst LANGAGES_PREFIX = \"Langages :\";\nconst VISUALISATION_PREFIX = \"Visualisation :\";\n\nlet langagesPara = null;\nlet visualisationPara = null;\n\nfor (const p of paragraphs.items) {\n  const t = (p.text || \"\").trim();\n  if (langagesPara === null && t.indexOf(LANGAGES_PREFIX) === 0) {\n    langagesPara = p;\n  }\n  if (visualisationPara === null && t.indexOf(VISUALISATION_PREFIX) === 0) {\n    visualisationPara = p;\n  }\n}\n\nif (!langagesPara || !visualisationPara) {\n  throw new Error(\n    \"Could not locate the 'Langages' and/or 'Visualisation' paragraphs.\"\n  );\n}\n\nconst langagesText = langagesPara.text;\n\n// Insert a fresh paragraph carrying the \"Langages\" text right before the\n// \"Visualisation\" paragraph (it inherits that paragraph's formatting, which\n// matches the original \"Langages\" paragraph's formatting in this document).\nvisualisationPara.insertParagraph(langagesText, \"Before\");\n\n// Remove the paragraph from its old location (right after \"MLOps\").\nlangagesPara.delete();\n\nawait context.sync();\n", "ps1": "# Move the \"Langages : ...\" skills line so it becomes the first item of the\n# \"COMPETENCES TECHNIQUES\" list (i.e. place it right before the\n# \"Visualisation : tableau\" paragraph) instead of sitting between \"MLOps\"\n# and \"Bases de donn\u00e9es\".\n$d = $word.ActiveDocument\n\n$langagesPrefix = \"Langages :\"\n$visualisationPrefix = \"Visualisation :\"\n\n$langagesIndex = -1\n$visualisationIndex = -1\n$langagesText = $null\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    $trimmed = $t.Trim()\n    if ($langagesIndex -eq -1 -and $trimmed.StartsWith($langagesPrefix)) {\n        $langagesIndex = $i\n        $langagesText = $trimmed\n    }\n    if ($visualisationIndex -eq -1 -and $trimmed.StartsWith($visualisationPrefix)) {\n        $visualisationIndex = $i\n    }\n}\n\nif ($langagesIndex -eq -1 -or $visualisationIndex -eq -1) {\n    throw \"Could not locate the 'Langages' and/or 'Visualisation' paragraphs.\"\n}\n\n# Insert a new paragraph right before \"Visualisation\" and give it the\n# \"Langages\" paragraph's text (the new paragraph inherits \"Visualisation\"'s\n# formatting, which matches the original \"Langages\" paragraph here).\n$visualisationRange = $d.Paragraphs.Item($visualisationIndex).Range\n$visualisationRange.InsertParagraphBefore()\n$d.Paragraphs.Item($visualisationIndex).Range.Text = $langagesText\n\n# Inserting a paragraph before \"Visualisation\" shifts every paragraph from\n# that point on (including the original \"Langages\" one, wherever it was)\n# down by one index; delete the original occurrence, paragraph mark\n# included, to leave a single \"Langages\" paragraph in its new spot.\nif ($langagesIndex -ge $visualisationIndex) {\n    $origLangagesIndex = $langagesIndex + 1\n} else {\n    $origLangagesIndex = $langagesIndex\n}\n$d.Paragraphs.Item($origLangagesIndex).Range.Delete()\n"}
